$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; this shifts existing rows 18-134 down to 19-135
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new data record
$ws.Cells.Item(18, 1).Value = 7
$ws.Cells.Item(18, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(18, 3).Value = "Ñuble"
$ws.Cells.Item(18, 4).Value = 45061
$ws.Cells.Item(18, 5).Value = 16
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100108
$ws.Cells.Item(18, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(18, 9).Value = 100108002
$ws.Cells.Item(18, 10).Value = "Mango"
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 50
$ws.Cells.Item(18, 14).Value = 8000
$ws.Cells.Item(18, 15).Value = 8000
$ws.Cells.Item(18, 16).Value = 8000
$ws.Cells.Item(18, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(18, 18).Value = "Perú"
$ws.Cells.Item(18, 19).Value = 2000
$ws.Cells.Item(18, 20).Value = 4
